# Append a new parameter/unit row to the bottom of the sheet,
# mirroring the existing "ScannedParameter" / "ScannedParameterUnit" list layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row right after the current data (row 53 -> 54)
$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "ZS_NI"
$ws.Cells.Item($newRow, 2).Value = "V"
